$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('B6').Value = '[''ar'' ''de'' ''en'']'
$ws.Range('D6').Value = '[''ar'' ''ca'' ''es'' ''id'']'
$ws.Range('E6').Value = '[''en'' ''ar'' ''undetected'' ''af'']'
$ws.Range('F6').Value = '[''ar'' ''undetected'']'
$ws.Range('G6').Value = '[''fr'' ''so'' ''ar'']'
$ws.Range('H6').Value = '[''en'' ''ar'' ''it'' ''af'' ''fr'']'
$ws.Range('I6').Value = '[''ar'' ''es'' ''en'' ''fr'' ''de'' ''id'' ''undetected'' ''ca'' ''tr'' ''pl'' ''tl'']'
$ws.Range('J6').Value = '[''en'' ''ar'' ''pt'']'
$ws.Range('K6').Value = '[''en'' ''ar'' ''ja'' ''id'' ''fr'']'
$ws.Range('L6').Value = '[''ar'' ''undetected'']'
$ws.Range('M6').Value = '[''en'' ''es'' ''ar'' ''it'' ''sv'']'
$ws.Range('N6').Value = '[''undetected'' ''ar'']'
$ws.Range('O6').Value = '[''ar'' ''en'' ''pl'' ''undetected'' ''id'']'
$ws.Range('P6').Value = '[''en'' ''undetected'']'
$ws.Range('Q6').Value = '[''ar'' ''es'' ''fr'' ''en'' ''undetected'']'
$ws.Range('R6').Value = '[''ar'' ''en'' ''fr'' ''undetected'' ''nl'' ''tr'' ''id'' ''ca'' ''de'' ''es'']'
$ws.Range('S6').Value = '[''ar'' ''undetected'' ''en'' ''ko'' ''ca'' ''id'']'
$ws.Range('T6').Value = '[''ar'' ''so'' ''en'' ''undetected'']'
$ws.Range('U6').Value = '[''ar'' ''es'' ''en'' ''undetected'' ''pt'' ''hu'' ''fa'' ''ca'' ''so'' ''et'' ''id'']'
$ws.Range('V6').Value = '[''es'' ''en'' ''ar'' ''undetected'' ''ca'' ''so'' ''tr'' ''de'' ''fr'']'
$ws.Range('W6').Value = '[''en'' ''ar'' ''de'' ''undetected'']'
$ws.Range('X6').Value = '[''ar'' ''es'' ''en'']'
$ws.Range('Y6').Value = '[''ar'' ''es'' ''undetected'' ''en'']'

$ws.Range('B7').Value = '[nan ''فعاليات'' ''شكرا'' ''احتفلت'' ''مبروك'' ''جانب'' ''اجواء'' ''الجماهير'' ''Die''
 ''تمنياتي'' ''Heute'' ''Today'' ''حضور'' ''World'' ''أقامت'' ''Thank'' ''اقامت'' ''سفارة''
 ''برنامج'' ''Frau'' ''Foreign'' ''وزارة'']'
$ws.Range('C7').Value = '[''رئيس'' ''بيان'' nan]'
$ws.Range('D7').Value = '[''La'' ''#qatar2022'' ''Así'' nan ''سفارة'' ''Hoy'' ''Viví'' ''▶️'' ''🎧'']'
$ws.Range('E7').Value = '[''As'' nan ''The'' ''Wishing'' ''During'' ''Hard'' ''Human'' ''And'' ''I'' ''We'' ''In''
 "It''s" ''بمناسبة'' ''Congratulations'' ''TODAY:'' ''What'' ''Hospitality''
 ''.@FIFAWorldCup'' ''Last'' ''Alongside'' ''On'' ''Football'' ''An'' ''Good'' ''Are''
 ''Today'' ''A'' ''Qatar'' ''#Qatar'' ''WATCH:'' "Qatar''s" ''This'' ''Go'' ''Finally!''
 ''Great'' ''His'' ''Unity'' ''⌛️5'' ''Look'' ''For'' ''Welcome'' ''Pleased'' ''Soccer''
 ''Second'' ''Empowered'' ''We’re'' ''Deputy'' ''“Qatar'' "HH''s" ''Moved'' ''#QFC''
 ''H.E.'' ''HE'']'
$ws.Range('F7').Value = '[nan ''رفعتم'' ''نفخر'' ''ألف'' ''انضمام'' ''أجواء'' ''حضرة'' ''أداء'' ''نشاط'' ''من'' ''قطر''
 ''سفارة'' ''اليوم'' ''نهنئ'' ''عندي'' ''مقتطف'' ''انطلقت'' ''انتظرونا.'' ''لم'' ''10''
 ''🔴بـــيــان'' ''شعار'' ''دولة'' ''رئيس'' ''الأمين'' ''يُعد'' ''أعلنت'' ''اللجنة'']'
$ws.Range('G7').Value = '[''"Notre'' nan ''"Je'' ''#ConférenceUkraine'' "L''ambassadeur" ''Un'' ''التقى''
 ''Nouvelle'' ''Entre'' ''🇶🇦'' ''أقام'']'
$ws.Range('H7').Value = '[nan ''Sheikh'' ''What'' ''The'' "Qatar''s" ''We'' ''Inspired'' ''In'' ''Deputy''
 ''International'' ''Spokesperson'' ''Committee'' ''Participated'' ''Pleased''
 ''Wishing'' ''Qatar,'' ''Today,'' ''بعد'' ''الحلم'' ''Great'' ''HE'' ''Statement'' ''FIFA''
 ''Welcome'' ''Delighted'' ''5'' ''Excited'' ''#QatarAirways'' ''As'' ''Infographic''
 ''To'' ''Son'' ''Qatar'' ''His'' ''Highlights'' ''On'' "HH''s" ''I'' ''whom'' ''Minister'']'
$ws.Range('I7').Value = '[nan ''Qatar'' ''Statement'' ''بيان'' ''رأيتُ'' ''Spokesperson'' ''المتحدث'' ''أقام''
 ''Deputy'' ''Minister'' ''مساعد'' ''HE'' ''لجنة'' ''Assistant'' ''The'' ''حالة''
 ''Secretary-General'' ''Disbursing'' ''انتخاب'' ''الأمين'' ''وزير'' ''رسالة'' ''وثمن''
 ''Gowning'' ''بيان|'' ''أتقدم'' ''دولة'' ''أعلنت'' ''سموّ'' ''شكرن'' ''البيان'' ''سمو''
 ''سفارة'' ''الأسطورة'' ''@majedalansari'' ''Iranian'' ''HH'' ''En'' ''كالعهد'' ''حضور''
 ''يأتي'' ''#وحدتنا_مصدر_قوتنا'' ''نائب'' ''📽️نائب'' ''سررت'' ''من'' ''We'' ''Pleased''
 ''#Our_Unity_Source_of_Our_Strength'' ''تنتهي'' ''بمناسبة'' ''This''
 ''#اليوم_الوطني_القطري'' ''Their'' ''قطر'' ''Félicitations'' ''Asistente''
 ''Viceprimer'' ''I'' ''التقيت'' ''Doha'' ''Lusail'' ''سعدت'' ''#GenderEquality''
 ''Delighted'' ''اختيار'' ''🎥|'' ''In'' ''Día'' ''Discurso'' ''Al-Kuwari:'' ''#MOFAQatar''
 ''Permanent'' "Qatar''s" ''President'' ''📽️Deputy'' ''شارك'' ''في'' ''كما'' ''THE''
 ''Qatar,'' ''يحتضن'' ''Through'' ''#الأمم_المتحدة'' ''International'' ''Remarks''
 ''PCOC'' ''الوفد'' ''تأتي'' ''اليوم'' ''المدير'' ''الكواري:'' ''Great'' ''UN'' ''الأمم''
 ''Estado'' ''(لئن'' ''أهنئ'' ''Portavoz'' ''Committee'' ''“The'' ''Conditions'' ''El''
 ''Briefed'' ''Esto'' ''#اليوم_الدولي_للتضامن_مع_الشعب_الفلسطيني'' ''تؤكد'' ''دعمت''
 ''تواصل'' ''Over'' ''#QNA_Video'' ''Brazilian'' ''🎥'' ''سعادة'' ''Participated'' ''Met''
 ''#الخارجية_القطرية'' ''ضمن'' ''Secretario'' ''As'' ''U.S.'' ''افتتحت'' ''Today,''
 ''نبارك'' ''Live'' ''يبدأ'' ''بث'' ''البث'' ''Regularity'' ''يعكس'' ''Fifth'' ''Zambia''
 ''Looking'' ''الحوار'' ''أتطلع'' ''"وجعلناكم'' ''بعد'' ''Foreign'' ''اللجنة'' ''المبعوث''
 ''دعم'' ''Gran'' ''نجدد'' ''أطيب'' ''✔️'' ''#HamadPort'' ''Embajada'' ''Embassy'' ''Desde''
 ''بقي'' ''#QNA_Infographic'' ''الإعلان'' ''#Doha'' ''#WorldCupQatar2022'' ''دشنا''
 ''@MBA_AlThani_'' ''🆕Bring'' ''#قنا_فيديو'' ''#قنا_إنفوجرافيك''
 ''#مونديال_قطر_2022..'' ''Be'' ''أحر'' ''بينما'' ''اللحظة'' ''Our'' ''Acaba'' ''بوصولك''
 ''Z'' ''الممثل'' ''Special'' ''خلال'' ''🔟'' ''Always'' ''المجموعة'' ''قبل'' ''مشاركة''
 ''#قطر'' ''After'' ''These'' ''وحدة'' ''تخرج'' ''Second'' ''Glad'' ''Congratulations''
 ''Merchandise'' ''#Football'' ''كرة'' ''Japanese'' ''Thank'' ''Infographic'' ''Major''
 ''Very'' ''📸|'' ''Read'' ''2/'' ''Empowered'' ''Did'' ''#FIFA'' ''Dans'' ''إنفوجراف'' ''To''
 ''لقراءة'' ''Son'' ''Qatar’s'' ''Lors'' ''On'' ''Pleasure'' "Aujourd''hui," ''During''
 ''.@ILOQatar'' ''Los'' ''Arab'' ''القادة'' ''Transforming'' ''According'' ''H.E'' ''ILO''
 ''منظمة'' ''كلمة'' ''نهنئ'' ''تستضيف'' ''Parte'' ''صورة'' ''Part'' ''جانب'' ''FIFA'' ''رئيس''
 ''Avec'' ''Ministerio'' ''وزارة'' ''موجز'' ''Today'' ''Inauguration'' ''Inauguración''
 ''افتتاح'' ''His'' ''التقى'' ''#MOCIQATAR,'' ''نظمت'' ''#MOCIQATAR''
 ''@QatarMission_Ge'' ''⏰'' ''❗'' ''H.E.'' ''#Qatar'' ''Mogadishu'' ''Vice-President''
 ''مدير'' ''جامعة'' ''سفارات'' ''Speaker'' ''وقعت'' ''Infograph'' ''Somali'' ''Ministry''
 ''Un'' ''State'' ''البدء'' ''انعقاد'' ''ممثل'' ''#UNESCO'' ''MOPH'' ''على'' ''أكد'' ''ملك''
 "HH''s" ''#DYK❓️'' ''تناول'' ''استضافة🇶🇦'' ''يؤكد'' ''مسؤولة'' ''توقيع'' ''Minostro''
 ''Moved'' ''نشارك'' ''Une'' ''دشنت'' ''الخارجية'' ''تقديرا'' ''بطلنا'' ''Happy'' ''حضرة''
 ''Chairman'' ''MOT'' ''⏳1'' ''Comandante'' ''⏳'' ''أقل'' ''Less'' ''Tanzanian'' ''Indian''
 ''للمرة'' ''رئيسة'' ''Tajikistan'' ''#قنا_انفوجرافيك'' "Iran''s" ''▶️'' ''◀️''
 "#Qatar''s" ''🏆'' ''Le'' ''✨'' ''✅'' ''تؤمن'' ''Sudanese'' ''يسعدُ'' ''#Qatar🇶🇦''
 ''Commercial'' ''La'' ''Secretary'' ''الرئيس'' ''الشيخة'' ''Association'' ''NHRC''
 ''Shura'' ''إنفوجراف|'' ''#مجلس_الشورى'' ''سفير'' ''Ministerial'' ''حضرت'' ''African''
 ''Education'' ''استلمت'' ''GANHRI'' ''Turkish'' ''مراسم'' ''#H.E.'' ''سعيد'' ''Russian''
 ''#FIFAWorldCup'' ''التصريحات'' ''Felicito'' ''جولة'' ''Lebanese'' ''Gulf'' ''At'' ''My''
 ''QNL'' ''Vicepresidente'' ''نائبة'' ''Vice'' ''#Shura_Council''
 ''#مكتبة_قطر_الوطنية'' ''Indonesian'' ''برنامج'' ''التقت'' ''Sharing'' ''Las''
 ''Executive'' ''Friendship'' ''السيد'' ''Join'' ''National'' ''وكيل'' ''يشارك'' ''¡No''
 ''Ministro'' ''القوات'' ''تعرّفوا'' ''It'' ''All'' ''Add'' ''فــخــر'' ''OUR'' ''Hamad''
 ''مشاورات'' ''Les'' ''#اليوم-الدولي_للمعلمين'' "#International_Teachers''_Day"
 ''وأغتنم'' ''Such'' ''paraguay'' ''فوز'' ''#Ashghal'' ''لقطات'' ''Museum'' ''Paraguay''
 ''الباراغواي'' ''Prime'']'
$ws.Range('J7').Value = '[''1.4'' ''سفارة'' ''This'' ''HH'' ''HE'' ''كلمة'' ''The'' ''We'' nan ''Inaugurado'' ''SE''
 ''Inaugurated'' ''#WorldCupQatar2022'' ''Labor'' ''#QNA_Video''
 ''#QNA_Infographic'' ''Group'' ''Discover'' ''To'' ''Qatar'']'
$ws.Range('K7').Value = '[''The'' ''تهنئكم'' nan ''أبرز''
 ''駐日カタール国大使館は、天皇陛下の63歳のお誕生日に際し、陛下並びに日本国民の皆様に心よりお祝い申し上げます。陛下の御健康と日本国民の皆様の益々のご繁栄を祈念申し上げます。''
 ''***تنــويه'' ''سفارة'' ''إنفوجراف|'' ''#قطر'' ''قطر'' ''الدكتور'' ''Qatar,''
 ''Exchange'' ''本日2月2日木曜日、カタール国と日本の間で、一般パスポート保持者に対する査証の相互免除に関する覚書が交わされました。''
 ''تبادل'' ''Japanese'' ''🎥|'' ''رئيس'' ''Met'' ''نائب'' ''سررت'' ''Statement'' ''بيان''
 ''بيان|'' ''موجز'' ''statement:'' ''سفارات'' ''اليوم'' ''H.E.''
 ''12月14日(水曜日)、ハッサン・ビン・モハメド・ラフィ・アルエマーディ駐日カタール国大使が帝国ホテルにてカタール建国記念式典を開催しました。日本政府の高官や各国の大使、ビジネスマン、カタール人留学生等、多くの人が参加しました。''
 ''أقـام'' ''His'' ''استقبل''
 ''12月15日（木曜日）、ハッサン・ビン・モハメド・ラフィ・アルエマーディ駐日カタール国大使は、経済産業省にて西村康稔経済産業大臣を表敬しました。会談では、両国間の関係の強化及び発展について話し合われ、双方の共通の関心事項についても述べられました。''
 ''سعادة'' ''カタール国建国記念レセプションが行われました🎊'' ''武井外務副大臣のカタール国ナショナルデー・レセプション出席'' ''In''
 ''International'' ''FIFAワールド杯'' ''Her''
 ''高円宮妃殿下はFIFAワールドカップ観戦のため訪問していたカタールで、日本代表が16強入りしたスペイン戦を応援されました。'' ''الأميرة''
 ''دعمت'' ''HE'' ''مرحباً'' ''Part'' ''Fifth''
 ''山梨県早川町の辻一幸町長が日本最古の宿「慶雲館」にて夕食会を開催し、ハッサン・ビン・モハメド・ラフィー・アルエマーディ駐日カタール国大使や湾岸諸国の大使らを招待しました。''
 ''حضر''
 ''長崎幸太郎山梨県知事がハッサン・ビン・モハメド・ラフィー・アルエマーディ駐日カタール国大使を含む湾岸諸国の大使らを山梨県へ招待し、米倉山太陽光発電所PR施設「ゆめソーラー館やまなし」を訪問しました。''
 ''ハッサン・ビン・モハメド・ラフィー・アルエマーディ駐日カタール国大使が長崎幸太郎''
 ''ハッサン・ビン・モハメド・ラフィー・アルエマーディ駐日カタール国大使が'' ''لقاء'' ''新しいバージョンをお届けします'' ''Mr.''
 ''本日11月11日、ハッサン・ビン・モハメド・ラフィー・アルエマーディ駐日カタール国大使が議員会館において、木原稔衆議院議員を表敬しました。両国関係の強化および発展について話し合われ、双方の共通の関心事についても述べられました。''
 ''未来をつなぐ湾岸都市''
 ''11月9日夜、ハッサン・ビン・モハメド・フィ・アルエマーディ駐日カタール国大使は、ワールドカップに出場する日本代表選手らを成田空港で見送りました。''
 ''想像を超える'' ''@MBA_AlThani_'' ''2022年11月4日、'' ''HRH'' ''سمو'' ''#FIFA'' ''Dans'' ''في''
 ''إستقبل'' ''شارك'' ''カタール国は豪政府によるエルサレム首都認定の撤回を歓迎'' ''دولة'' ''التصريحات'' ''مراسم'']'
$ws.Range('L7').Value = '[''سفارات'' nan ''بعدسة'' ''سمو'' ''#وحدتنا_مصدر_قوتنا'' ''جياني'' ''رئيس''
 ''انفانتينو:'' ''إشادات'' ''#قنا_انفوجرافيك'' ''اليوم'' ''البيان'' ''مجلة''
 ''#قطر_للطاقة'' ''سيتم'' ''وزير'' ''دولة'' ''🎥|'' ''#كأس_العالم_قطر_2022'' ''ليلة''
 ''دعم'' ''#قنا_فيديو'' ''نجدد'' ''لطالما'' ''#كأس_العالم_2022'' ''#قنا_إنفوجرافيك''
 ''جمعية'' ''تم'' ''إنفوجراف'' ''🎥'' ''المتحدث'' ''في'' ''عدد'']'
$ws.Range('M7').Value = '[''Qatar'' ''“Once'' ''Statement'' nan ''Qatar’s'' ''What'' ''ARGENTINA'' ''We'' ''The''
 ''Welcome'' ''📽️Deputy'' ''This'' ''#Qatar'' ''Congratulations'' ''History'' ''𝐓𝐡𝐞''
 ''Club'' ''Almost'' ''A'' ''Did'' ''In'' ''QFFD'']'
$ws.Range('N7').Value = '[nan]'
$ws.Range('O7').Value = '[nan ''Warm'' ''Spokesperson'' ''المتحدث'' ''بعد'' ''يتعين'' ''تحت'' ''نهنئ'' ''I'' ''Z''
 ''أهنئ'' ''We'' ''@MBA_AlThani_'' ''سعادة'' ''HE'' ''Read'' ''إنفوجراف'' ''Infographic''
 ''حضر'' ''Qatar,'' ''مشاورات'']'
$ws.Range('P7').Value = '[nan ''ARGENTINA'' ''We'' ''𝐇𝐈𝐒𝐓𝐎𝐑𝐘'' ''The'' ''Happy'' ''History'' ''Congratulations''
 ''Well'' ''Great'' ''As'' ''In'' ''Welcome'' ''HE'' ''#Qatar'' ''Deputy'' ''H.E.'']'
$ws.Range('Q7').Value = '[''المتحدث'' ''بيان'' nan ''تحولت'' ''بيان|'' ''دولة'' ''أتقدم'' ''نائب'' ''وزير'' ''#شاهد''
 ''الأمين'' ''انتخاب'' ''سفارات'' ''بأجواء'' ''بمناسبة'' ''تنتهي'' ''أهنئ'' ''حضور''
 ''جانب'' ''سررت'' ''التقيت'' ''📽️نائب'' ''مساعد'' ''سعدت'' ''اختيار'' ''🎥|'' ''الكواري:''
 ''المدير'' ''بمشاركة'' ''اليوم'' ''تأتي'' ''نفخر'' ''الوفد'' ''#الخارجية_القطرية''
 ''لجنة'' ''اشتراطات'' ''فوز'' ''كلمة'' ''تواصل'' ''دعمت'' ''🎥'' ''تفاعل'' ''سعادة''
 ''احتفالات'' ''قطر'' ''افتتحت'' ''بث'' ''البث'' ''يعكس'' ''مشاركة'' ''الحوار'' ''سفارة''
 ''أتطلع'' ''بالتعاون'' ''أقامت'' ''زامبيا'' ''Parte'' ''Dans'' ''Part'' ''انطلق'' ''حضرة''
 ''بعد'' ''#عالوعد'' ''المبعوث'' ''مرحباً'' ''دعم'' ''أطيب'' ''معلومات'' ''سمو'' ''كما''
 ''إعلان'' ''أضفنا'' ''الجمهور'' ''قائمة'' ''لحظة'' ''يشكل'' ''دشنا'' ''عندما'' ''كرة''
 ''@MBA_AlThani_'' ''القادة'' ''إن'' ''منظمة'' ''تستضيف'' ''موجز'' ''اللجنة'' ''افتتاح''
 ''تعرّفوا'' ''مدير'' ''باعتبارها'' ''لا'' ''إنفوجراف'' ''وزارة'' ''أكد'' ''يؤكد'' ''في''
 ''نشارك'' ''تناول'' ''من'' ''دشنت'' ''الخارجية'' ''قِصةُ'' ''قصة'' ''أبرز''
 ''#وحدتنا_مصدر_قوتنا'' ''رئيس'' ''تؤمن'' ''حضرت'' ''سفير'' ''جولة'' ''شهد'']'
$ws.Range('R7').Value = '[nan ''تمكنت'' ''Through'' ''The'' ''تنتهي'' ''أهنئ'' ''Qatar'' ''أتقدم'' ''نتقدم''
 ''سفارة'' ''Al'' ''نائب'' ''Deputy'' ''📽️Deputy'' ''Le'' ''من'' ''Inspired'' ''رفعتم''
 ''كفو'' ''المدير'' ''دولة'' ''🎥|'' ''Pleased'' ''سعدت'' ''نفخر'' ''سمو'' ''بذلت'' ''Met''
 ''استوحي'' ''As'' ''سعادة'' ''Participated'' ''المتحدث'' ''حضرة'' ''بعد'' ''"وجعلناكم''
 ''أتطلع'' ''Looking'' ''جانب'' ''مشاركة'' ''📽️نائب'' ''باقٍ'' ''أضفنا'' ''بوصولك'' ''يشكل''
 ''Inaugurated'' ''أحر'' ''I'' ''بيان'' ''Statement'' ''شارك'' ''دشنا'' ''HE''
 ''@MBA_AlThani_'' ''ministre'' ''افتتاح'' ''It'' ''Son'' ''#MOFAQatar''
 ''#الخارجية_القطرية'' ''لقراءة'' ''To'' ''@MofaQatar_EN'' ''Dans'' ''Infographic''
 ''إنفوجراف'' ''سررت'' ''🎥'' ''Minister'' ''وزير'' ''We'' ''باعتبارها'' ''HH'' ''في'' ''قطر''
 ''His'' ''In'' ''نشارك'' ''Today,'' ''دشنت'' ''HE.'' ''شكراً'' ''Felicito'' ''حضرت'' ''H.E''
 ''إنفوجراف|'' ''مراسم'' ''التصريحات'' ''لقطات'' ''Pleasant'']'
$ws.Range('S7').Value = '[nan ''سفارة'' ''احتفال'' ''Al'' ''During'' ''خيمتنا'' ''نائب'' ''قامت'' ''주한카타르대사관은''
 ''Spokesperson'' ''As'' ''Minister'' ''Deputy'' ''Qatar'' ''دولة'' ''بعد'' "La''eeb"
 ''#World_Cup_Qatar_2022'' ''#Qatar2022,'' ''دعم'' ''8️⃣'' ''That'' ''دشنا''
 ''@MBA_AlThani_'' ''Ahead'' ''Marvel'' ''#وحدتنا_مصدر_قوتنا'']'
$ws.Range('T7').Value = '[nan ''🇦🇷Amb'' ''وفد'' ''أهنئ'' ''انتخاب'' ''Qatar'' ''سررت'' ''بعد'' ''أطيب'' ''أتقدم''
 ''سمو'' ''سعدت'' ''Might'' ''لاهاي،'']'
$ws.Range('U7').Value = '[''سفارة'' nan ''#Our_Unity_Source_of_Our_Strength'' ''Día'' ''اقامت'' ''La''
 ''#Qatar2022'' ''El'' ''Experience'' ''Metro?'' ''Viceprimer'' ''En'' ''Al-Kuwari:''
 ''Discurso'' ''Asistente'' ''Estado'' ''Travelling'' ''Looking'' ''ℹ️'' ''UNA'' ''¡Así''
 ''🇫🇷'' ''Lo'' ''Buenas'' ''USG'' ''SE'' ''عبد'' ''Visit'' ''Portavoz'' ''discurso'' ''Esto''
 ''¡Solo'' ''Clasificados'' ''¡Estos'' ''¡Inglaterra'' ''Cuando'' ''"Creo'' ''Explore''
 ''From'' ''Energy'' ''Una'' ''¡Súmate'' ''📍'' ''Meeetro?'' ''¡Esto'' ''📸'' ''¡Se''
 ''¡Momento'' ''¡Empate'' ''¡Victoria'' ''⚽️'' ''1,'' ''Así'' ''¡Oficialmente'' ''0''
 "It''s" ''Gran'' ''Feel'' ''Con'' ''Qatar'' ''Su'' ''Transporte'' ''Competition'' ''📲''
 ''5'' ''💾'' ''كلمة'' ''Posts'' ''Visitors'' ''Sin'' ''Los'' ''If'' ''Acaba'' ''Touchdown''
 ''Ministerio'' ''Un'' ''Si'' ''🏟'' ''Falta'' ''¿Cómo'' ''Inauguración'' ''Here'' ''✨'' ''6''
 ''🇶🇦'' ''Have'' ''❗'' ''Msheireb,'' ''⏰'' ''Did'' ''¡FALTAN'' ''¡5️⃣'' ''¡FALTA'' ''¿Vas''
 ''¿Qué'' ''32'' ''¿Todavía'' ''¡Fiesta'' ''Solo'' ''¡El'' ''🏆'' ''اكتشف'' "Surf''s"
 ''¿Quién'' ''Dos'' ''Nueva'' ''Take'' ''Only'' ''👇'' ''🚨'' ''Felicito'' ''When'' ''🔸''
 ''@MIAQatar'' "Qatar''s" ''Happy'' ''Save'' ''¡No'']'
$ws.Range('V7').Value = '[''Qatar'' ''Statement'' ''بيان'' nan ''بيان|'' ''نائب'' ''Deputy'' ''انتخاب''
 ''#Our_Unity_Source_of_Our_Strength'' ''En'' ''#اليوم_الوطني_القطري'' ''بمناسبة''
 ''يأتي'' ''#وحدتنا_مصدر_قوتنا'' ''احتفال'' ''تنتهي'' ''La'' ''احتفلت'' ''Logo''
 ''Viceprimer'' ''📽️نائب'' ''سررت'' ''I'' ''In'' "Qatar''s" ''Pleased''
 ''Secretary-General'' ''Doha'' ''اختيار'' ''الأمين'' ''دولة'' ''🎥|'' ''سعدت'' ''Día''
 ''Asistente'' ''Discurso'' ''Al-Kuwari:'' ''📽️Deputy'' ''THE'' ''Qatar,'' ''PCOC''
 ''Minister'' ''Remarks'' ''Assistant'' ''البيان'' ''في'' ''الكواري:'' ''المدير'' ''من''
 ''مساعد'' ''قطر'' ''الأمم'' ''UN'' ''نفخر'' ''Estado'' ''Felicito'' ''أتقدم'' ''سمو'' ''El''
 ''discurso'' ''Portavoz'' ''Committee'' ''Spokesperson'' ''The'' ''Over'' ''U.S.'' ''As''
 ''Secretario'' ''Gran'' ''HH'' ''Su'' ''#QNA_Video'' ''Las'' ''Ibrahim''
 ''#WorldCupQatar2022'' ''Wael'' ''Ahmed'' ''#مونديال_قطر_2022'' ''السفير'' ''بوصولك''
 ''Acaba'' ''During'' ''الوفد'' ''Los'' ''Foreign'' ''Arab'' ''Special'' ''ILO'' ''On''
 ''تستضيف'' ''منظمة'' ''المبعوث'' ''المتحدث'' ''القادة'' ''FALTAN'' ''Parte''
 ''Ministerio'' ''GCC'' ''FIFA'' ''Part'' ''وزارة'' ''اللجنة'' ''🎥'' ''Un'' ''جانب''
 ''Inauguración'' ''وزير'' ''افتتاح'' ''¡Estas'' ''No'' ''سفارات'' ''مدير'' ''أكد'' "HH''s"
 ''Today,'' ''يؤكد'' ''تناول'' ''نشارك'' ''H.E'' ''باقي'' ''30'' ''حضرة'' ''Comandante''
 ''Tajikistan'' ''أهنئ'' ''¡No'' ''🏆'' ''Infographic'' ''Ministry'' ''تؤمن'' ''سفارة''
 ''إنفوجراف|'' ''Esta'' ''Ministerial'' ''African'' ''حضرت'' ''انعقاد'' ''المجموعة''
 ''#قنا_انفوجرافيك'' ''مراسم'' ''التصريحات'' ''Les'' ''Una'' ''paraguay'' ''Desde''
 ''Paraguay'' ''Museum'' ''#Qatar'' ''لقطات'' ''📌'']'
$ws.Range('W7').Value = '[''Statement'' ''Spokesperson'' ''Assistant'' ''Deputy'' ''Minister'' ''Qatar'' nan
 ''تنتهي'' ''Pleased'' ''بمناسبة'' ''#Our_Unity_Source_of_Our_Strength'' ''I''
 ''Doha'' ''📽️Deputy'' ''In'' ''PCOC'' ''Qatar,'' ''THE'' ''The'' ''International''
 ''Committee'' ''Over'' ''Participated'' ''As'' ''U.S.'' ''Today,'' ''Was'' ''HH'' ''Great''
 ''بعد'' ''HE'' ''Be'' ''Japanese'' ''Read'' ''Special'' ''Qatar’s'' ''Infographic''
 ''Arab'' ''ILO'' ''On'' ''FIFA'' ''Inauguration'' ''Here'' "HH''s" ''Part'' ''Ministry''
 ''H.E'' ''Ministerial'' ''Come'' ''🎥'' ''Paraguay'']'
$ws.Range('X7').Value = '[nan ''Statement'' ''El'']'
$ws.Range('Y7').Value = '[nan ''#قنا_فيديو'' ''As'']'
